$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("overview_testcases")
$ws1.Activate()
$p = $excel.ActiveWindow.Panes.Item(1)
Write-Host "pane:" $p
Write-Host "type:" $p.GetType()
$tlc = $p.TopLeftCell
Write-Host "tlc:" $tlc
